# Updated cryptos list - refresh price & volume(1h) figures, and swap the
# FraxShare / PaxDollar rows (42 <-> 43) to reflect the new ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price/Volume cells in this sheet are stored as text. Excel's COM
# layer auto-converts plain-decimal-looking strings (e.g. "211.50") into
# numbers when assigned directly, so for those we briefly force a Text
# number format, assign the value, then restore the default ("Normal")
# style so no stray formatting is left behind.
function Set-TextValue($addr, $val) {
    if ($val -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = $val
        $ws.Range($addr).Style = "Normal"
    } else {
        $ws.Range($addr).Value = $val
    }
}

function Set-Row($Row, $D, $E) {
    if ($D) {
        Set-TextValue "D$Row" $D
    }
    if ($E) {
        Set-TextValue "E$Row" $E
    }
}

Set-Row 2  "27.512.33"  "  +1.67%  "
Set-Row 3  "1.573.82"   "  +0.43%  "
Set-Row 4  $null        "  -1.54%  "
Set-Row 5  "211.50"     "  +1.46%  "
Set-Row 6  $null        "  +0.16%  "
Set-Row 7  $null        "  -1.45%  "
Set-Row 8  "22.87"      "  +3.42%  "
Set-Row 9  "0.250"      "  +0.61%  "
Set-Row 10 $null        "  +0.23%  "
Set-Row 11 "0.0870"     "  +1.36%  "
Set-Row 12 "1.797.67"   "  +0.40%  "
Set-Row 13 "1.569.04"   "  +0.15%  "
Set-Row 14 $null        "  -0.48%  "
Set-Row 15 "0.519"      "  -0.10%  "
Set-Row 16 "27.499.77"  "  +1.61%  "
Set-Row 17 "62.47"      "  +0.96%  "
Set-Row 18 "226.61"     "  +5.13%  "
Set-Row 19 "7.52"       "  +1.22%  "
Set-Row 20 $null        "  +0.17%  "
Set-Row 21 "0.991"      "  -1.50%  "
Set-Row 22 $null        "  -0.66%  "
Set-Row 23 $null        "  +2.39%  "
Set-Row 24 $null        "  +0.57%  "
Set-Row 25 "150.43"     "  -2.32%  "
Set-Row 26 "15.18"      "  +0.93%  "
Set-Row 27 $null        "  -0.33%  "
Set-Row 28 $null        "  +1.77%  "
Set-Row 29 $null        "  -1.46%  "
Set-Row 30 $null        "  +0.92%  "
Set-Row 31 $null        "  -0.33%  "
Set-Row 32 $null        "  +0.46%  "
Set-Row 33 "1.455.62"   "  +2.29%  "
Set-Row 34 "3.14"       "  -1.87%  "
Set-Row 35 $null        "  +3.32%  "
Set-Row 36 $null        "  -0.13%  "
Set-Row 37 $null        "  -0.86%  "
Set-Row 38 $null        "  +0.39%  "
Set-Row 39 $null        "  +1.30%  "
Set-Row 40 $null        "  +0.11%  "
Set-Row 41 "2.37"       "  -0.62%  "

# Rows 42 & 43 swap content: FraxShare moves up to row 42, PaxDollar moves
# down to row 43, each with refreshed price/volume figures.
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D42" "5.65"
Set-TextValue "E42" "  -2.96%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D43" "0.991"
Set-TextValue "E43" "  -1.41%  "

Set-Row 44 $null        "  +6.91%  "
Set-Row 45 "0.976"      "  -3.01%  "
Set-Row 46 "64.35"      "  -0.51%  "
Set-Row 47 "1.709.72"   "  +0.36%  "
Set-Row 48 $null        "  +0.25%  "
Set-Row 49 "0.0₆0105"   "  +0.93%  "
Set-Row 50 $null        "  +1.45%  "
Set-Row 51 "0.0946"     "  -1.78%  "
